$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.405.08"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.509.69"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "312.61"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "98.22"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").Value = "35.00"
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "7.14"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "2.894.26"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "2.504.32"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "15.12"
$ws.Range("E16").Value = "  -5.24%  "
$ws.Range("D17").Value = "0.802"
$ws.Range("E17").Value = "  -4.01%  "
$ws.Range("D18").Value = "42.374.69"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "6.55"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("D20").Value = "0.0₃0930"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "11.96"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("D22").Value = "68.31"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "239.78"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("D27").Value = "25.17"
$ws.Range("E27").Value = "  -5.04%  "
$ws.Range("E28").Value = "  -5.00%  "
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "37.47"
$ws.Range("E30").Value = "  -7.63%  "
$ws.Range("D31").Value = "156.32"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "5.78"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Value = "2.64"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("D35").Value = "0.0777"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D38").Value = "17.44"
$ws.Range("E38").Value = "  -4.88%  "
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").Value = "21.31"
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("D45").Value = "1.996.09"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  -4.47%  "
$ws.Range("D47").Value = "8.85"
$ws.Range("D48").Value = "2.745.91"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("D49").Value = "78.23"
$ws.Range("E49").Value = "  -3.79%  "
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").Value = "70.95"
$ws.Range("E51").Value = "  -3.27%  "
